# Work Profile and new tenant support
# Appends new sprint-run history rows to the "AMSIN" and "AMS" worksheets,
# and normalizes the formatting of the previously-last row on "AMSIN"
# (row 10) to match the rest of the data rows.
#
# Pattern used for every cell: set the VALUE first (forcing a text
# NumberFormat beforehand for text columns so date-looking strings like
# "2023-02-20" are not auto-converted into date serials), and only
# afterwards paste-special the FORMATS from an existing "template" cell
# that already carries the desired style. Because the value is written
# before the formats-only paste, the paste cannot clobber it.

$wb = $excel.ActiveWorkbook

function Write-HistoryRow {
    param(
        $ws,
        [int]$rowNum,
        [int]$templateRow,
        $dtWs,
        [int]$dtTemplateRow,
        $a, $b, $c, $d, $e, $f, $g
    )

    # --- Column A (Run Date) : literal text ---
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $a
    $ws.Range("A" + $templateRow).Copy()
    $cellA.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- Column B (Run Time) : numeric date-time serial ---
    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.Value = $b
    $dtWs.Range("B" + $dtTemplateRow).Copy()
    $cellB.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- Column C (Sprint Name) : literal text ---
    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $c
    $ws.Range("C" + $templateRow).Copy()
    $cellC.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- Columns D-G (numeric) ---
    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.Value = $d
    $ws.Range("D" + $templateRow).Copy()
    $cellD.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.Value = $e
    $ws.Range("E" + $templateRow).Copy()
    $cellE.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $cellF = $ws.Cells.Item($rowNum, 6)
    $cellF.Value = $f
    $ws.Range("F" + $templateRow).Copy()
    $cellF.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $cellG = $ws.Cells.Item($rowNum, 7)
    $cellG.Value = $g
    $ws.Range("G" + $templateRow).Copy()
    $cellG.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# =======================================================================
# Sheet "AMSIN" : style-normalize row 10, append rows 11-18
# =======================================================================
$wsAMSIN = $wb.Worksheets.Item("AMSIN")

# Row 10 already holds data; only its formatting changes (pick up the
# same general style already used by rows 8-9), plus a tiny correction
# to the Run Time value. Re-use the row-writer so every column goes
# through the same value-then-format sequence.
Write-HistoryRow $wsAMSIN 10 9 $wsAMSIN 9 "2023-02-20" 44977.42567164352 "173aadharflow" 34 34 0 1.35

Write-HistoryRow $wsAMSIN 11 9 $wsAMSIN 9 "2023-03-09" 44994.56723940972 "aadhar174fstcycle" 34 34 0  1.46
Write-HistoryRow $wsAMSIN 12 9 $wsAMSIN 9 "2023-03-13" 44998.44610030093 "174aadharflow"     34 34 0  1.43
Write-HistoryRow $wsAMSIN 13 9 $wsAMSIN 9 "2023-03-30" 45015.69776930555 "175scnadhar"        34 32 2  2.39
Write-HistoryRow $wsAMSIN 14 9 $wsAMSIN 9 "2023-03-31" 45016.50551984954 "175fnlaad"          34 33 1  1.65
Write-HistoryRow $wsAMSIN 15 9 $wsAMSIN 9 "2023-04-06" 45022.65158956018 "176newaadha"        33 31 2  3.57
Write-HistoryRow $wsAMSIN 16 9 $wsAMSIN 9 "2023-04-07" 45023.67214146991 "176fstadh"          33 32 1  1.57
Write-HistoryRow $wsAMSIN 17 9 $wsAMSIN 9 "2023-04-07" 45023.67613211805 "176fstadh"          33 23 10 0.8

# Row 18: plain/default style on A,C,D,E,F,G (no template applied there),
# only the Run Time column keeps the datetime style.
$cellA18 = $wsAMSIN.Cells.Item(18, 1)
$cellA18.NumberFormat = "@"
$cellA18.Value = "2023-04-07"

$cellB18 = $wsAMSIN.Cells.Item(18, 2)
$cellB18.Value = 45023.68523743405
$wsAMSIN.Range("B9").Copy()
$cellB18.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAMSIN.Cells.Item(18, 3).NumberFormat = "@"
$wsAMSIN.Cells.Item(18, 3).Value = "176fstadh"
$wsAMSIN.Cells.Item(18, 4).Value = 33
$wsAMSIN.Cells.Item(18, 5).Value = 24
$wsAMSIN.Cells.Item(18, 6).Value = 9
$wsAMSIN.Cells.Item(18, 7).Value = 1.05

# =======================================================================
# Sheet "AMS" : append rows 8-11
# =======================================================================
$wsAMS = $wb.Worksheets.Item("AMS")

Write-HistoryRow $wsAMS 8  7 $wsAMSIN 9 "2023-03-13" 44998.56252199074 "174betaaadh" 34 33 1 1.39
Write-HistoryRow $wsAMS 9  7 $wsAMSIN 9 "2023-03-13" 44998.85407626157 "173liveaadd" 34 33 1 1.19
Write-HistoryRow $wsAMS 10 7 $wsAMSIN 9 "2023-03-31" 45016.56328677083 "175btaadh"   34 33 1 1.23
Write-HistoryRow $wsAMS 11 7 $wsAMSIN 9 "2023-03-31" 45016.8308946875  "175aaddev"   33 32 1 1.19

Write-Output "Edit complete"
